# Apply "Emociones emitidas por agentes" update to frames_view.xlsx
# Adds 6 new rows (8-13) describing the "whatsapp" segment, each with a
# hyperlink in column C, mirroring the existing afiliaciones rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row data: segmento | categoria | url
$rows = @(
    @{ Row = 8;  Cat = "3d";                         Url = "https://projector.tensorflow.org/?config=https://gist.githubusercontent.com/luiska1803/e65e57239606051a55acecc9753d29d6/raw/d7940b8bcfa124c1f7de469da063bc735e5bd092/Wapp%2520Cafam" },
    @{ Row = 9;  Cat = "agente_experiencia_emitida";  Url = "https://htmlpreview.github.io/?https://github.com/DavidIngx/html/blob/master/cf_entregable/Whatsapp/Palabras_emociones_user.html" },
    @{ Row = 10; Cat = "sentimientos_usuarios";       Url = "https://htmlpreview.github.io/?https://github.com/DavidIngx/html/blob/master/cf_entregable/Whatsapp/Experiencia%20General%20User.html" },
    @{ Row = 11; Cat = "palabras_clave";              Url = "https://htmlpreview.github.io/?https://github.com/DavidIngx/html/blob/master/cf_entregable/Whatsapp/Key%20words%20User.html" },
    @{ Row = 12; Cat = "2d topics";                   Url = "https://htmlpreview.github.io/?https://github.com/DavidIngx/html/blob/master/cf_entregable/Whatsapp/vis_LDA_Chat_wp_usuario.html" },
    @{ Row = 13; Cat = "n_gramas";                    Url = "https://htmlpreview.github.io/?https://github.com/DavidIngx/html/blob/master/cf_entregable/Whatsapp/(N-GRAMAS)_%20User.html" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Range("A$rowNum").Value = "whatsapp"
    $ws.Range("B$rowNum").Value = $r.Cat
    $ws.Range("C$rowNum").Value = $r.Url

    $ws.Hyperlinks.Add($ws.Range("C$rowNum"), $r.Url, [Type]::Missing, [Type]::Missing, $r.Url)
    $ws.Range("C$rowNum").Style = "Hipervínculo"
}

